$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.201916098594666
$ws.Range("B1").Value = 2.611344575881958
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.175978422164917
$ws.Range("E1").Value = 1.170597195625305
